$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "Meta description" paragraph (near the top) and the
# "Please create an image..." paragraph (the very last paragraph)
# by content, rather than assuming fixed indices.
# ------------------------------------------------------------------
$metaIndex = -1
$imagePromptIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($metaIndex -eq -1 -and $txt.StartsWith("Meta description")) {
        $metaIndex = $i
    }
    if ($txt.StartsWith("Please create an image")) {
        $imagePromptIndex = $i
    }
}

$newHeadingText = "Play Drago: Jewels of Fortune for Free - Expert Slot Game Review"
$newBodyText = "Read our expert review of Drago: Jewels of Fortune and play for free. Experience the dragon-themed design, streak respin feature, and free spins with multipliers."
$boldLabel = "Meta description"

# ------------------------------------------------------------------
# Step 1: Capture the "Meta description" run's bold-only formatting
# (just the "<w:b/>" run, not the plain run after it) so the new
# heading paragraph we build further down reuses that exact run
# formatting instead of whatever happens to sit at the insertion
# point (which would otherwise bleed in unwanted italics / styles).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item($metaIndex)
$metaParaStart = $metaPara.Range.Start
$boldRange = $d.Range($metaParaStart, $metaParaStart + $boldLabel.Length)
$boldFormatted = $boldRange.FormattedText

# ------------------------------------------------------------------
# Step 2: Insert a brand-new paragraph right after the paragraph that
# precedes the image-prompt paragraph (i.e. right before it), and
# stamp it with the captured bold formatting. Resetting Style to
# "Normal" first keeps it from inheriting the preceding paragraph's
# list-bullet style while still letting FormattedText bring in the
# character-level (bold) formatting cleanly.
# ------------------------------------------------------------------
$beforeImagePrompt = $d.Paragraphs.Item($imagePromptIndex - 1)
$beforeImagePrompt.Range.InsertParagraphAfter()

$newCount = $d.Paragraphs.Count
$newHeadingPara = $d.Paragraphs.Item($imagePromptIndex)
$newHeadingPara.Style = "Normal"
$newHeadingRange = $newHeadingPara.Range
$newHeadingRange.FormattedText = $boldFormatted

$newHeadingTextRange = $d.Range($newHeadingRange.Start, $newHeadingRange.Start + $boldLabel.Length)
$newHeadingTextRange.Text = $newHeadingText

# ------------------------------------------------------------------
# Step 3: The image-prompt paragraph has now been pushed one slot
# further down; replace its text with the former meta-description
# body copy while keeping its own (italic) run formatting intact.
# ------------------------------------------------------------------
$imagePromptIndex = $imagePromptIndex + 1
$imagePromptPara = $d.Paragraphs.Item($imagePromptIndex)
$imagePromptRange = $imagePromptPara.Range
$imagePromptReplaceRange = $d.Range($imagePromptRange.Start, $imagePromptRange.End)
$imagePromptReplaceRange.Text = $newBodyText

# ------------------------------------------------------------------
# Step 4: Remove the original "Meta description" paragraph from the
# top of the document now that its content has been relocated.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item($metaIndex)
$metaPara.Range.Delete()
